# Update the workbook "Översikt HANINGE" per the recorded diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The "Förändrad" (changed) date in column C changes from 45172 to 45175
#    for every data row (rows 2 through 91).
for ($r = 2; $r -le 91; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45175
}

# 2) A new case "A 34417-2023" replaces the old row 12 entry and becomes the
#    new row 11 (with updated figures / an extra species and no "Markägare").
#    The former row 11 ("A 32532-2020") is pushed down to row 12 (only its
#    "Förändrad" date changes, already handled above).

# --- New row 11: A 34417-2023 ---
$ws.Range("A11").Value2 = "A 34417-2023"
$ws.Range("B11").Value2 = 45139
$ws.Range("C11").Value2 = 45175
$ws.Range("D11").Value2 = "STOCKHOLMS LÄN"
$ws.Range("E11").Value2 = "HANINGE"
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value2 = 3.8
$ws.Range("H11").Value2 = 1
$ws.Range("I11").Value2 = 3
$ws.Range("J11").Value2 = 1
$ws.Range("K11").Value2 = 1
$ws.Range("L11").Value2 = 0
$ws.Range("M11").Value2 = 0
$ws.Range("N11").Value2 = 0
$ws.Range("O11").Value2 = 2
$ws.Range("P11").Value2 = 1
$ws.Range("Q11").Value2 = 6
$ws.Range("R11").Value2 = "Porslinsblå spindling`r`nBarrviolspindling`r`nDropptaggsvamp`r`nKornknutmossa`r`nSårläka`r`nBlåsippa"
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 34417-2023.xlsx")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 34417-2023.png")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 34417-2023.docx")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 34417-2023.docx")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 34417-2023.docx")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 34417-2023.docx")'

# --- New row 12: A 32532-2020 (previously row 11, shifted down) ---
$ws.Range("A12").Value2 = "A 32532-2020"
$ws.Range("B12").Value2 = 44018
$ws.Range("C12").Value2 = 45175
$ws.Range("D12").Value2 = "STOCKHOLMS LÄN"
$ws.Range("E12").Value2 = "HANINGE"
$ws.Range("F12").Value2 = "Kommuner"
$ws.Range("G12").Value2 = 2.9
$ws.Range("H12").Value2 = 3
$ws.Range("I12").Value2 = 2
$ws.Range("J12").Value2 = 2
$ws.Range("K12").Value2 = 0
$ws.Range("L12").Value2 = 0
$ws.Range("M12").Value2 = 0
$ws.Range("N12").Value2 = 0
$ws.Range("O12").Value2 = 2
$ws.Range("P12").Value2 = 0
$ws.Range("Q12").Value2 = 5
$ws.Range("R12").Value2 = "Spillkråka`r`nTalltita`r`nGranbarkgnagare`r`nVedticka`r`nRevlummer"
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 32532-2020.xlsx")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 32532-2020.png")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 32532-2020.docx")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 32532-2020.docx")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 32532-2020.docx")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 32532-2020.docx")'
